$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227, shifting existing rows 227-327 down to 228-328
$ws.Range("A227").EntireRow.Insert()

# Populate the newly inserted row 227 with the new data record
$ws.Cells.Item(227, 1).Value = 4
$ws.Cells.Item(227, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(227, 3).Value = "Los Lagos"
$ws.Cells.Item(227, 4).Value = 44992
$ws.Cells.Item(227, 5).Value = 10
$ws.Cells.Item(227, 6).Value = 100112039
$ws.Cells.Item(227, 7).Value = "Ciboulette"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 240
$ws.Cells.Item(227, 11).Value = 3500
$ws.Cells.Item(227, 12).Value = 3500
$ws.Cells.Item(227, 13).Value = 3500
$ws.Cells.Item(227, 14).Value = "`$/docena de atados"
$ws.Cells.Item(227, 15).Value = "Región Metropolitana"
$ws.Cells.Item(227, 16).Value = 1167
$ws.Cells.Item(227, 17).Value = 3
$ws.Cells.Item(227, 18).Value = "Hortaliza"
